$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 5-7 (cases C4, C5, C6) get the same numeric formatting ("0.0000E+00")
# already used by rows 2-4 (cases C1-C3).
$ws.Range("B5:G7").NumberFormat = $ws.Range("B4:G4").NumberFormat

# Row 5 -> Case C4
$ws.Range("B5").Value = 20411500
$ws.Range("C5").Value = 111389000
$ws.Range("D5").Value = 24995500
$ws.Range("E5").Value = 2.5749
$ws.Range("F5").Value = 0.57780699999999996
$ws.Range("G5").Value = 4.31691

# Row 6 -> Case C5
$ws.Range("B6").Value = 20411500
$ws.Range("C6").Value = 22277700
$ws.Range("D6").Value = 999822
$ws.Range("E6").Value = 0.51497999999999999
$ws.Range("F6").Value = 0.023112299999999999
$ws.Range("G6").Value = 4.31691

# Row 7 -> Case C6
$ws.Range("B7").Value = 20411500
$ws.Range("C7").Value = 222777000
$ws.Range("D7").Value = 99982200
$ws.Range("E7").Value = 5.1497999999999999
$ws.Range("F7").Value = 2.3112300000000001
$ws.Range("G7").Value = 4.31691

$wb.Save()
